# Commit message: "Added 1 more country to exluded"
#
# The row for Timor-Leste (ISO3 "TLS") is removed from the
# iso3CountryCoordinates sheet (it was row 170: TLS, -8.874217, 125.727539,
# Timor-Leste). Deleting the whole row shifts every following row up by
# one and Excel automatically drops the now-unused "TLS"/"Timor-Leste"
# shared-string entries on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the Timor-Leste row (A170:D170) first, matching how a user would
# highlight the row before deleting it, then delete the entire row so
# everything below shifts up.
$row = $ws.Range("A170:D170")
$row.Select()
$row.EntireRow.Delete()
